$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new blank column before column N (shifts N->O, O->P, P->Q)
$ws.Range("N1").EntireColumn.Insert()

# Update the selection to reflect where the user clicked after the edit
$ws.Range("T6").Select()
